$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.907.33"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.549.67"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.58"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.24"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0587"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.770.76"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.549.76"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.74"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "26.902.85"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.70"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.75"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.25"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.26"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "1.417.26"
$ws.Range("E33").Value = "  +3.33%  "
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.71"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.43"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "1.684.40"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.55"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0517"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  -0.01%  "
